# Update statistics values on the "SummaryReport" sheet (Maximum / Std. Deviation / 90 Percent
# columns) to reflect refreshed load-test results.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SummaryReport")

$ws.Range("E2").Value = 31.687
$ws.Range("F2").Value = 6.466
$ws.Range("G2").Value = 31.415

$ws.Range("E3").Value = 0.128
$ws.Range("G3").Value = 0.078

$ws.Range("E4").Value = 0.114
$ws.Range("F4").Value = 0.015
$ws.Range("G4").Value = 0.086

$ws.Range("E5").Value = 0.733
$ws.Range("F5").Value = 0.048
$ws.Range("G5").Value = 0.72

$ws.Range("E6").Value = 0.843
$ws.Range("F6").Value = 0.06
$ws.Range("G6").Value = 0.759

$ws.Range("E7").Value = 0.163
$ws.Range("F7").Value = 0.016
$ws.Range("G7").Value = 0.144

$ws.Range("E8").Value = 0.168
$ws.Range("F8").Value = 0.018
$ws.Range("G8").Value = 0.128

$ws.Range("E9").Value = 0.215
$ws.Range("F9").Value = 0.021
$ws.Range("G9").Value = 0.15

$ws.Range("E10").Value = 0.15
$ws.Range("F10").Value = 0.017
$ws.Range("G10").Value = 0.109

$ws.Range("F11").Value = 0.01
$ws.Range("G11").Value = 0.077

$ws.Range("E12").Value = 0.183
$ws.Range("F12").Value = 0.015
$ws.Range("G12").Value = 0.155

$ws.Range("E13").Value = 0.19
$ws.Range("F13").Value = 0.018
$ws.Range("G13").Value = 0.13

$ws.Range("E14").Value = 0.147
$ws.Range("G14").Value = 0.141

$ws.Range("E15").Value = 10.407
$ws.Range("F15").Value = 0.041
$ws.Range("G15").Value = 10.375

$ws.Range("E16").Value = 25.728
$ws.Range("F16").Value = 0.043
$ws.Range("G16").Value = 25.702

$ws.Range("E17").Value = 31.687
$ws.Range("F17").Value = 0.08
$ws.Range("G17").Value = 31.471

$ws.Range("E18").Value = 21.361
$ws.Range("F18").Value = 0.091
$ws.Range("G18").Value = 21.325

$ws.Range("E19").Value = 21.868
$ws.Range("F19").Value = 0.082
$ws.Range("G19").Value = 21.815

$ws.Range("E20").Value = 16.673
$ws.Range("F20").Value = 0.019
